$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark attendance for 7th March (row 8) as present (TRUE) for all people
$ws.Range("B8:F8").Value = $true

# Update the active selection to F8 to match the edited cell
$ws.Range("F8").Select()
